# Update crypto price/volume columns (D, E) per the latest data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces text storage so numeric-looking prices (e.g. trailing
# zeros like 295.20) keep their exact printed form instead of being parsed as numbers.

$ws.Range("D2").Value = "'39.327.94"
$ws.Range("E2").Value = "  -2.27%  "
$ws.Range("D3").Value = "'2.193.69"
$ws.Range("E3").Value = "  -6.56%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'295.20"
$ws.Range("E5").Value = "  -4.50%  "
$ws.Range("D6").Value = "'81.67"
$ws.Range("E6").Value = "  -4.73%  "
$ws.Range("D7").Value = "'0.510"
$ws.Range("E7").Value = "  -3.85%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  -4.31%  "
$ws.Range("D10").Value = "'0.0768"
$ws.Range("E10").Value = "  -6.79%  "
$ws.Range("D11").Value = "'29.02"
$ws.Range("E11").Value = "  -4.88%  "
$ws.Range("D12").Value = "'47.05"
$ws.Range("E12").Value = "  -10.88%  "
$ws.Range("E13").Value = "  -2.64%  "
$ws.Range("D14").Value = "'2.538.38"
$ws.Range("E14").Value = "  -6.03%  "
$ws.Range("D15").Value = "'6.22"
$ws.Range("E15").Value = "  -3.57%  "
$ws.Range("D16").Value = "'13.92"
$ws.Range("E16").Value = "  -6.58%  "
$ws.Range("D17").Value = "'2.196.85"
$ws.Range("E17").Value = "  -6.31%  "
$ws.Range("D18").Value = "'0.709"
$ws.Range("E18").Value = "  -6.30%  "
$ws.Range("D19").Value = "'39.208.41"
$ws.Range("E19").Value = "  -2.35%  "
$ws.Range("D20").Value = "'0.0₃0867"
$ws.Range("E20").Value = "  -4.62%  "
$ws.Range("E21").Value = "  -6.85%  "
$ws.Range("D22").Value = "'64.66"
$ws.Range("E22").Value = "  -4.88%  "
$ws.Range("D23").Value = "'10.25"
$ws.Range("E23").Value = "  -5.52%  "
$ws.Range("D24").Value = "'224.73"
$ws.Range("E24").Value = "  -4.84%  "
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("D26").Value = "'2.39"
$ws.Range("E26").Value = "  -6.67%  "
$ws.Range("E27").Value = "  -0.98%  "
$ws.Range("D28").Value = "'22.52"
$ws.Range("E28").Value = "  -4.88%  "
$ws.Range("E29").Value = "  -2.15%  "
$ws.Range("D30").Value = "'9.06"
$ws.Range("E30").Value = "  -2.60%  "
$ws.Range("D31").Value = "'148.42"
$ws.Range("E31").Value = "  -2.33%  "
$ws.Range("D32").Value = "'31.77"
$ws.Range("E32").Value = "  -9.65%  "
$ws.Range("D33").Value = "'0.999"
$ws.Range("E33").Value = "  -0.18%  "
$ws.Range("D34").Value = "'4.79"
$ws.Range("E34").Value = "  -7.40%  "
$ws.Range("D35").Value = "'0.0692"
$ws.Range("E35").Value = "  -4.93%  "
$ws.Range("E36").Value = "  -5.01%  "
$ws.Range("E37").Value = "  -3.57%  "
$ws.Range("D38").Value = "'15.31"
$ws.Range("E38").Value = "  -4.37%  "
$ws.Range("D39").Value = "'0.0957"
$ws.Range("E39").Value = "  -4.99%  "
$ws.Range("D40").Value = "'2.62"
$ws.Range("E40").Value = "  -6.43%  "
$ws.Range("E41").Value = "  -4.79%  "
$ws.Range("E42").Value = "  -6.44%  "
$ws.Range("D43").Value = "'1.897.99"
$ws.Range("E43").Value = "  -2.76%  "
$ws.Range("E44").Value = "  -9.28%  "
$ws.Range("D45").Value = "'0.0258"
$ws.Range("E45").Value = "  -4.05%  "
$ws.Range("D46").Value = "'16.03"
$ws.Range("E46").Value = "  -10.44%  "
$ws.Range("E47").Value = "  -4.72%  "
$ws.Range("E48").Value = "  -4.43%  "
$ws.Range("D49").Value = "'71.41"
$ws.Range("E49").Value = "  -0.26%  "
$ws.Range("D50").Value = "'2.406.99"
$ws.Range("E50").Value = "  -5.92%  "
$ws.Range("D51").Value = "'86.98"
$ws.Range("E51").Value = "  -6.94%  "
